$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark that currently sits between
#    the "c" and "ontainer" runs of the "Sequence container" paragraph.
#    It will be re-created later at the end of the new "XML Task" item.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Find the "Sequence container" paragraph (Control Flow Tasks > Containers)
#    and insert the new list of Control Flow Tasks after it.
# ------------------------------------------------------------------
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Sequence container") {
        $anchor = $p
    }
}

if ($anchor -eq $null) {
    throw "Could not locate 'Sequence container' paragraph"
}

$items = @(
    @{ Text = "Tasks"; Level = 1 },
    @{ Text = "Active X Script Task"; Level = 2 },
    @{ Text = "Analysis Services Execute DDL Task"; Level = 2 },
    @{ Text = "Analysis Services Processing Task"; Level = 2 },
    @{ Text = "Bulk Insert Task"; Level = 2 },
    @{ Text = "Data Flow Task"; Level = 2 },
    @{ Text = "Data Mining Query Task"; Level = 2 },
    @{ Text = "Data Profiling Task"; Level = 2 },
    @{ Text = "Execute DTS 2000 Package Task"; Level = 2 },
    @{ Text = "Execute Package Task"; Level = 2 },
    @{ Text = "Execute Process Task"; Level = 2 },
    @{ Text = "Execute SQL Task"; Level = 2 },
    @{ Text = "File System Task"; Level = 2 },
    @{ Text = "FTP Task"; Level = 2 },
    @{ Text = "Message Queue Task"; Level = 2 },
    @{ Text = "Script Task"; Level = 2 },
    @{ Text = "Send Mail Task"; Level = 2 },
    @{ Text = "Transfer Database Task"; Level = 2 },
    @{ Text = "Transfer Error Messages Task"; Level = 2 },
    @{ Text = "Transfer Jobs Task"; Level = 2 },
    @{ Text = "Transfer Logins Task"; Level = 2 },
    @{ Text = "Transfer Master Stored Procedures Task"; Level = 2 },
    @{ Text = "Transfer SQL Server Objects Task"; Level = 2 },
    @{ Text = "Web Service Task"; Level = 2 },
    @{ Text = "WMI Data Reader Task"; Level = 2 },
    @{ Text = "WMI Event Watcher Task"; Level = 2 },
    @{ Text = "XML Task"; Level = 2 }
)

$prev = $anchor
$lastNewPara = $null
foreach ($item in $items) {
    $prev.Range.InsertParagraphAfter()
    $idx = $prev.Index + 1
    $newPara = $d.Paragraphs.Item($idx)
    $newPara.Style = "ListParagraph"
    $newPara.Range.ListFormat.ListLevelNumber = $item.Level
    $newPara.Range.Text = $item.Text
    $prev = $newPara
    $lastNewPara = $newPara
}

# Final empty bullet paragraph (ilvl 0) after "XML Task"
$prev.Range.InsertParagraphAfter()
$idx = $prev.Index + 1
$emptyPara = $d.Paragraphs.Item($idx)
$emptyPara.Style = "ListParagraph"
$emptyPara.Range.ListFormat.ListLevelNumber = 1

# ------------------------------------------------------------------
# 3. Re-create the "_GoBack" bookmark, collapsed, right after the text
#    of the new "XML Task" paragraph (i.e. at the end of lastNewPara).
#    NOTE: the runtime mishandles a collapsed bookmark placed exactly
#    one character before a paragraph end, so a sentinel character is
#    temporarily appended, the bookmark is added just before it, and
#    the sentinel is removed again.
# ------------------------------------------------------------------
$sentinel = [char]1
$lastNewPara.Range.InsertAfter($sentinel)
$bmPos = $lastNewPara.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
$sentinelRange = $d.Range($bmPos, $bmPos + 1)
$sentinelRange.Delete()
